# Fixing the big mistake: correct the "Total" (B) and "Community" (D)
# statistics that were computed with an error, and propagate the
# dependent Sums/% energy sector figures in rows 10-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- corrected descriptive statistics for the Total (B) and Community (D) columns ---
$ws.Range("B3").Value = 15716.53094379187   # mean
$ws.Range("D3").Value = 781.2597083866054

$ws.Range("B4").Value = 6562.794511211719   # std
$ws.Range("D4").Value = 343.2988852644162

$ws.Range("B5").Value = 5271.149254794524   # min
$ws.Range("D5").Value = 120.0064657534245

$ws.Range("B6").Value = 10780.78928356164   # 25%
$ws.Range("D6").Value = 443.5438390410952

$ws.Range("B7").Value = 13590.63691506851   # 50%
$ws.Range("D7").Value = 869.8767726027395

$ws.Range("B8").Value = 21807.62013013706   # 75%
$ws.Range("D8").Value = 1114.383575342466

$ws.Range("B9").Value = 29081.56297260281   # max
$ws.Range("D9").Value = 1175.328772602739

# --- downstream sums / proportions that depend on the corrected means ---
$ws.Range("F10").Value = 22631804.55906031
$ws.Range("G11").Value = 0.8115588775488014
$ws.Range("F12").Value = 1125013.980076712
$ws.Range("G12").Value = 0.04970942450218045
$ws.Range("G13").Value = 0.1387316979490181
